# Updates cryptocurrency price (D) and 1h-volume-change (E) text values
# for the symbol list, mirroring the upstream GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'259.29"
$ws.Range("E2").Value = "'0.57%"
$ws.Range("D3").Value = "'26.94"
$ws.Range("E3").Value = "'-1.89%"
$ws.Range("D4").Value = "'4.676"
$ws.Range("E4").Value = "'2.33%"
$ws.Range("D5").Value = "'0.06022"
$ws.Range("E5").Value = "'2.29%"
$ws.Range("D6").Value = "'6.665"
$ws.Range("E6").Value = "'0.48%"
$ws.Range("D7").Value = "'0.8608"
$ws.Range("E7").Value = "'0.20%"
$ws.Range("D8").Value = "'0.9300"
$ws.Range("E8").Value = "'0.52%"
$ws.Range("D9").Value = "'0.1395"
$ws.Range("E9").Value = "'-1.13%"
$ws.Range("D10").Value = "'0.04948"
$ws.Range("E10").Value = "'34.47%"
$ws.Range("D11").Value = "'0.07010"
$ws.Range("E11").Value = "'-0.99%"
$ws.Range("D12").Value = "'0.03130"
$ws.Range("E12").Value = "'-1.35%"
$ws.Range("D13").Value = "'0.09140"
$ws.Range("E13").Value = "'-0.42%"
$ws.Range("D14").Value = "'0.001537"
$ws.Range("E14").Value = "'-0.02%"
$ws.Range("D15").Value = "'0.0006040"
$ws.Range("E15").Value = "'-0.14%"
$ws.Range("D16").Value = "'0.006120"
$ws.Range("E16").Value = "'0.49%"
$ws.Range("E17").Value = "'-1.51%"
$ws.Range("D18").Value = "'3.163"
$ws.Range("E18").Value = "'-1.28%"
$ws.Range("D20").Value = "'0.3112"
$ws.Range("E20").Value = "'0.21%"
$ws.Range("E21").Value = "'1.55%"
$ws.Range("D22").Value = "'4.135"
$ws.Range("E22").Value = "'7.00%"
$ws.Range("D23").Value = "'0.04240"
$ws.Range("E23").Value = "'0.65%"
$ws.Range("D24").Value = "'0.001215"
$ws.Range("E24").Value = "'-0.52%"
$ws.Range("D25").Value = "'0.004039"
$ws.Range("E25").Value = "'-6.10%"
$ws.Range("D26").Value = "'0.0001198"
$ws.Range("E26").Value = "'-0.12%"
$ws.Range("E27").Value = "'13.55%"
$ws.Range("D40").Value = "'0.03842"
$ws.Range("E40").Value = "'0.10%"
$ws.Range("E41").Value = "'1.17%"
$ws.Range("D42").Value = "'0.003885"
$ws.Range("E42").Value = "'-38.10%"
$ws.Range("D43").Value = "'0.002416"
$ws.Range("E43").Value = "'9.86%"
$ws.Range("D44").Value = "'0.01526"
$ws.Range("E44").Value = "'28.99%"
$ws.Range("D45").Value = "'0.00005102"
$ws.Range("E45").Value = "'-6.50%"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("E47").Value = "'-16.72%"
$ws.Range("D48").Value = "'0.1503"
$ws.Range("E48").Value = "'16.06%"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.02%"
